$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 (duplicate of row 2)
$ws.Range("A4").Value2 = 45644
$ws.Range("A4").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B4").Value2 = 2016
$ws.Range("C4").Value = 'Driver drowsiness detection using forehead brain signals'
$ws.Range("D4").Value = 'Journal of Engineering and Applied Science       '
$ws.Range("E4").Value2 = 63
$ws.Range("F4").Value = 'Publisher, https://www.scopus.com/inward/record.uri?eid=2-s2.0-85051085110&partnerID=40&md5=8b0c84f1b9b90f287f2d6ee0910cbdaf'
$ws.Range("G4").Value = ""
$ws.Range("H4").Value = 'abstract: Driver drowsiness is a major problem causing car accidents. Many studies assessed drowsiness indicators at different driving setups. Brain signals recorded at different head sites resulted high accuracies. Forehead sites facilitates practical implementation of drowsiness detection system. The detection of drowsiness after normal working day with normal sleep habits was the scope. In this study, subjects carried out simulated monotonous driving task for two hours after 6 PM. Brain signals and facial changes were recorded simultaneously. Recorded data for all subjects were segmented to five seconds epochs. Each epoch of brain signals was labelled as alert or drowsy in accordance to subject''s facial symptoms in video records. Features were extracted out of brain signals recorded at forehead sites (Fp1 and Fp2 with Fpz reference) by discrete wavelet decomposition. Support vector machine was used for classification. The highest classification accuracy was 84.73% from combined group of features consisting of, energy of Alpha, energy of Beta, standard deviations of detail coefficients at fifth level (Alpha) and fourth level (Beta). Support vector machine classifier setup was radial basis function kernel. Many reasons contribute in the variance of accuracy, the driving task setup, preparation of subject for driving, vigilance states labeling, and signal preprocessing. © 2016 Medwell Journals. All rights reserved.:       '
$ws.Range("I4").Value = 'C:\Users\rpb\OneDrive - ums.edu.my\research_related\0 eeg_trend_till24\eeg_review\Driver_Drowsiness_Detection_Using_Deep_Learning.pdf'
$ws.Range("J4").Value = ""
$ws.Range("K4").Value = 'Eldeib_A'
$ws.Range("L4").Value = 'Eldeib_A_2016'
$ws.Range("M4").Value = 'relevance'
$ws.Range("N4").Value = ""
$ws.Range("O4").Value = 'C:\Users\balan\IdeaProjects\academic_paper_maker\ragpdf\test.pdf'

# Row 5 (duplicate of row 3)
$ws.Range("A5").Value = 'A. Al-Ani and M. Mesbah'
$ws.Range("B5").Value2 = 2018
$ws.Range("C5").Value = 'EEG rhythm/channel selection for fuzzy rule-based alertness state characterization'
$ws.Range("D5").Value = 'Neural Computing and Applications       '
$ws.Range("E5").Value2 = 30
$ws.Range("F5").Value = 'Publisher'
$ws.Range("G5").Value = '10.1007/s00521-016-2835-1'
$ws.Range("H5").Value = 'abstract: The aim of the paper is to automatically select the optimal EEG rhythm/channel combinations capable of classifying human alertness states. Four alertness states were considered, namely ‘engaged’, ‘calm’, ‘drowsy’ and ‘asleep’. The features used in the automatic selection are the energies associated with the conventional rhythms, δ, θ, α, β and γ, extracted from overlapping windows of the different EEG channels. The selection process consists of two stages. In the first stage, the optimal brain regions, represented by sets of EEG channels, are selected using a simple search technique based on support vector machine (SVM), extreme learning machine (ELM) and LDA classifiers. In the second stage, a fuzzy rule-based alertness classification system (FRBACS) is used to identify, from the previously selected EEG channels, the optimal features and their supports. The IF–THEN rules used in FRBACS are constructed using a novel differential evolution-based search algorithm particularly designed for this task. Each alertness state is represented by a set of IF–THEN rules whose antecedent parts contain EEG rhythm/channel combination. The selected spatio-frequency features were found to be good indicators of the different alertness states, as judged by the classification performance of the FRBACS that was found to be comparable to those of the SVM, ELM and LDA classifiers. Moreover, the proposed classification system has the advantage of revealing simple and easy to interpret decision rules associated with each of the alertness states. © 2016, The Natural Computing Applications Forum.:       '
$ws.Range("I5").Value = ""
$ws.Range("J5").Value = 'https://www.scopus.com/inward/record.uri?eid=2-s2.0-85007492413&doi=10.1007%2fs00521-016-2835-1&partnerID=40&md5=e471b46f7066da5b9465b33b0d118386'
$ws.Range("K5").Value = 'Mesbah_A'
$ws.Range("L5").Value = 'Mesbah_A_2018'
$ws.Range("M5").Value = 'relevance'
$ws.Range("N5").Value = ""
$ws.Range("O5").Value = 'C:\Users\balan\IdeaProjects\academic_paper_maker\ragpdf\test2.pdf'

# Row 6 (duplicate of row 2)
$ws.Range("A6").Value2 = 45644
$ws.Range("A6").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B6").Value2 = 2016
$ws.Range("C6").Value = 'Driver drowsiness detection using forehead brain signals'
$ws.Range("D6").Value = 'Journal of Engineering and Applied Science       '
$ws.Range("E6").Value2 = 63
$ws.Range("F6").Value = 'Publisher, https://www.scopus.com/inward/record.uri?eid=2-s2.0-85051085110&partnerID=40&md5=8b0c84f1b9b90f287f2d6ee0910cbdaf'
$ws.Range("G6").Value = ""
$ws.Range("H6").Value = 'abstract: Driver drowsiness is a major problem causing car accidents. Many studies assessed drowsiness indicators at different driving setups. Brain signals recorded at different head sites resulted high accuracies. Forehead sites facilitates practical implementation of drowsiness detection system. The detection of drowsiness after normal working day with normal sleep habits was the scope. In this study, subjects carried out simulated monotonous driving task for two hours after 6 PM. Brain signals and facial changes were recorded simultaneously. Recorded data for all subjects were segmented to five seconds epochs. Each epoch of brain signals was labelled as alert or drowsy in accordance to subject''s facial symptoms in video records. Features were extracted out of brain signals recorded at forehead sites (Fp1 and Fp2 with Fpz reference) by discrete wavelet decomposition. Support vector machine was used for classification. The highest classification accuracy was 84.73% from combined group of features consisting of, energy of Alpha, energy of Beta, standard deviations of detail coefficients at fifth level (Alpha) and fourth level (Beta). Support vector machine classifier setup was radial basis function kernel. Many reasons contribute in the variance of accuracy, the driving task setup, preparation of subject for driving, vigilance states labeling, and signal preprocessing. © 2016 Medwell Journals. All rights reserved.:       '
$ws.Range("I6").Value = 'C:\Users\rpb\OneDrive - ums.edu.my\research_related\0 eeg_trend_till24\eeg_review\Driver_Drowsiness_Detection_Using_Deep_Learning.pdf'
$ws.Range("J6").Value = ""
$ws.Range("K6").Value = 'Eldeib_A'
$ws.Range("L6").Value = 'Eldeib_A_2016'
$ws.Range("M6").Value = 'relevance'
$ws.Range("N6").Value = ""
$ws.Range("O6").Value = 'C:\Users\balan\IdeaProjects\academic_paper_maker\ragpdf\test.pdf'

# Row 7 (duplicate of row 3)
$ws.Range("A7").Value = 'A. Al-Ani and M. Mesbah'
$ws.Range("B7").Value2 = 2018
$ws.Range("C7").Value = 'EEG rhythm/channel selection for fuzzy rule-based alertness state characterization'
$ws.Range("D7").Value = 'Neural Computing and Applications       '
$ws.Range("E7").Value2 = 30
$ws.Range("F7").Value = 'Publisher'
$ws.Range("G7").Value = '10.1007/s00521-016-2835-1'
$ws.Range("H7").Value = 'abstract: The aim of the paper is to automatically select the optimal EEG rhythm/channel combinations capable of classifying human alertness states. Four alertness states were considered, namely ‘engaged’, ‘calm’, ‘drowsy’ and ‘asleep’. The features used in the automatic selection are the energies associated with the conventional rhythms, δ, θ, α, β and γ, extracted from overlapping windows of the different EEG channels. The selection process consists of two stages. In the first stage, the optimal brain regions, represented by sets of EEG channels, are selected using a simple search technique based on support vector machine (SVM), extreme learning machine (ELM) and LDA classifiers. In the second stage, a fuzzy rule-based alertness classification system (FRBACS) is used to identify, from the previously selected EEG channels, the optimal features and their supports. The IF–THEN rules used in FRBACS are constructed using a novel differential evolution-based search algorithm particularly designed for this task. Each alertness state is represented by a set of IF–THEN rules whose antecedent parts contain EEG rhythm/channel combination. The selected spatio-frequency features were found to be good indicators of the different alertness states, as judged by the classification performance of the FRBACS that was found to be comparable to those of the SVM, ELM and LDA classifiers. Moreover, the proposed classification system has the advantage of revealing simple and easy to interpret decision rules associated with each of the alertness states. © 2016, The Natural Computing Applications Forum.:       '
$ws.Range("I7").Value = ""
$ws.Range("J7").Value = 'https://www.scopus.com/inward/record.uri?eid=2-s2.0-85007492413&doi=10.1007%2fs00521-016-2835-1&partnerID=40&md5=e471b46f7066da5b9465b33b0d118386'
$ws.Range("K7").Value = 'Mesbah_A'
$ws.Range("L7").Value = 'Mesbah_A_2018'
$ws.Range("M7").Value = 'relevance'
$ws.Range("N7").Value = ""
$ws.Range("O7").Value = 'C:\Users\balan\IdeaProjects\academic_paper_maker\ragpdf\test2.pdf'
